# PlayerPerformance_3745.xlsx edit script
# - Inserts a new "Player Info" sheet before "ODI Batting"
# - Renames MATCH_CARD_LINK -> MATCH_CODE and converts the URL column to
#   bare numeric match codes on both "ODI Batting" and "ODI Bowling"
# - Drops now-redundant empty INNING_NUMBER cells on "ODI Batting"
# - Appends a new "ODI Batting Extra" sheet with additional per-match stats

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, placed before "ODI Batting"
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($battingSheet)

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3745"
$playerInfo.Range("B2").Value = "Ranasinghe Arachchige Suranga Lakmal"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code,
#    remove empty INNING_NUMBER cells
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingRows = $battingSheet.Cells.Item($battingSheet.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value
    if ($link) {
        $pos = $link.LastIndexOf("=")
        if ($pos -ge 0) {
            $code = $link.Substring($pos + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }

    $inningCell = $battingSheet.Cells.Item($r, 2)
    if (-not $inningCell.Value) {
        $inningCell.Value = $null
    }
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingRows = $bowlingSheet.Cells.Item($bowlingSheet.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value
    if ($link) {
        $pos = $link.LastIndexOf("=")
        if ($pos -ge 0) {
            $code = $link.Substring($pos + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $code
        }
    }
}

# ---------------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet appended at the end
# ---------------------------------------------------------------------------
$extra = $wb.Worksheets.Add()
$extra.Name = "ODI Batting Extra"
$extra.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
$extra.Range("A1:F1").Font.Bold = $true

$extraData = @(
    @("4084", 10, "3", "0", "13.29%", "NO"),
    @("4096", 9, $null, $null, $null, "YES"),
    @("4098", $null, $null, $null, $null, "NO"),
    @("4099", 10, "0", "0", "0.47%", "NO"),
    @("4112", 10, "0", "0", "1.80%", "NO"),
    @("4116", 10, "0", "0", "0.64%", "NO"),
    @("4119", 10, $null, $null, $null, "NO"),
    @("4122", 8, $null, $null, $null, "YES"),
    @("4124", 10, "0", "0", "0.90%", "NO"),
    @("4182", 9, "1", "0", "2.59%", "NO"),
    @("4183", $null, $null, $null, $null, "NO"),
    @("4186", 9, "1", "0", "4.21%", "NO"),
    @("4187", 10, $null, $null, $null, "NO"),
    @("4188", $null, $null, $null, $null, "NO"),
    @("4193", 9, "1", "1", "16.13%", "NO"),
    @("4302", $null, $null, $null, $null, "NO"),
    @("4305", 10, "1", "0", "5.15%", "NO"),
    @("4309", 9, "2", "0", "7.46%", "NO"),
    @("4339", 10, "0", "0", "2.46%", "NO"),
    @("4451", 9, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraData) {
    $extra.Cells.Item($r, 1).NumberFormat = "@"
    $extra.Cells.Item($r, 1).Value = $row[0]

    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]

    $r++
}
